# Target change: split the run
#   " (if we want to see the request being displayed in the result page) and then click on "
# into three runs, changing the word "request" -> "response":
#   " (if we want to see the "  +  "response"  +  " being displayed in the result page) and then click on "
#
# A plain Find/Replace (or a plain Range.Text assignment) keeps everything inside a single
# run, because the host re-coalesces adjacent runs that share identical formatting when it
# saves. To make the new "response" text land in its own run (matching the target XML,
# where none of the three runs carry any rPr), we first stake out the "request" range with
# a temporary bookmark. Anchoring a bookmark to that sub-range forces the engine to keep it
# as a distinct run boundary once the text is overwritten, and removing the bookmark
# afterwards leaves no trace in the saved document.

$d = $word.ActiveDocument

$search = $d.Content
$search.Find.Execute("request being displayed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $search.Start
$target = $d.Range($start, $start + 7)   # "request"

$d.Bookmarks.Add("tmp_request_word", $target)
$target.Text = "response"
$d.Bookmarks("tmp_request_word").Delete()
